$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in row 12 with the new "Exp 16" experiment parameters
$ws.Range("A12").Value = "Exp 16"
$ws.Range("B12").Value = 0.6
$ws.Range("C12").Value = 1
$ws.Range("F12").Value = "Exp 16.png"

# Update the active selection to match the saved workbook state
$ws.Range("F13").Select()
